$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 618.94116
$ws.Range("I33").Value = 288.35715
$ws.Range("K33").Value = 288.35715
$ws.Range("M33").Value = -59.35714999999999

$ws.Range("H99").Value = 4309.4
$ws.Range("I99").Value = 540
$ws.Range("K99").Value = 1620
$ws.Range("M99").Value = -122

$ws.Range("H111").Value = 1029
$ws.Range("I111").Value = 1029
$ws.Range("K111").Value = 3087
$ws.Range("M111").Value = -20

$ws.Range("H132").Value = 2473.4912
$ws.Range("I132").Value = 2166.463
$ws.Range("K132").Value = 6499.389000000001
$ws.Range("M132").Value = -3969.389000000001

$ws.Range("H137").Value = 2893.6667
$ws.Range("I137").Value = 2893.6667
$ws.Range("K137").Value = 8681.000100000001
$ws.Range("M137").Value = -6131.000100000001

$ws.Range("H138").Value = 3058.1355
$ws.Range("I138").Value = 1584.0952
$ws.Range("J138").Value = 3872.7368
$ws.Range("K138").Value = 4752.2856
$ws.Range("L138").Value = 11618.2104
$ws.Range("M138").Value = 387.7143999999998
$ws.Range("N138").Value = -21898.2104

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2204
$ws.Range("I2").Value = 1823.6818
$ws.Range("K2").Value = 1823.6818
$ws.Range("M2").Value = -1710.6818

$ws.Range("H45").Value = 3627.238
$ws.Range("I45").Value = 2695.3333
$ws.Range("K45").Value = 2695.3333
$ws.Range("M45").Value = -2318.3333

$ws.Range("H61").Value = 5404.25
$ws.Range("J61").Value = 10565.857
$ws.Range("L61").Value = 10565.857
$ws.Range("N61").Value = -10989.857

$ws.Range("H63").Value = 7484.7144
$ws.Range("I63").Value = 3247
$ws.Range("J63").Value = 9179.8
$ws.Range("K63").Value = 3247
$ws.Range("L63").Value = 9179.8
$ws.Range("M63").Value = -2561
$ws.Range("N63").Value = -10551.8

$ws.Range("H66").Value = 7484.7144
$ws.Range("I66").Value = 3247
$ws.Range("J66").Value = 9179.8
$ws.Range("K66").Value = 16235
$ws.Range("L66").Value = 45899
$ws.Range("M66").Value = -12803
$ws.Range("N66").Value = -52763

$ws.Range("H110").Value = 2064.3
$ws.Range("I110").Value = 1830.375
$ws.Range("J110").Value = 3000
$ws.Range("K110").Value = 1830.375
$ws.Range("L110").Value = 3000
$ws.Range("M110").Value = 214.625
$ws.Range("N110").Value = -7090

$ws.Range("H116").Value = 2204
$ws.Range("I116").Value = 1823.6818
$ws.Range("K116").Value = 1823.6818
$ws.Range("M116").Value = 470.3181999999999

$ws.Range("H136").Value = 5404.25
$ws.Range("J136").Value = 10565.857
$ws.Range("L136").Value = 31697.571
$ws.Range("N136").Value = -36797.571

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2204
$ws.Range("I3").Value = 1823.6818
$ws.Range("K3").Value = 1823.6818
$ws.Range("M3").Value = -1709.6818

$ws.Range("H86").Value = 2241.0588
$ws.Range("I86").Value = 1999.909
$ws.Range("J86").Value = 2683.1667
$ws.Range("K86").Value = 1999.909
$ws.Range("L86").Value = 2683.1667
$ws.Range("M86").Value = -876.9090000000001
$ws.Range("N86").Value = -4929.1667

$ws.Range("H89").Value = 2241.0588
$ws.Range("I89").Value = 1999.909
$ws.Range("J89").Value = 2683.1667
$ws.Range("K89").Value = 9999.545
$ws.Range("L89").Value = 13415.8335
$ws.Range("M89").Value = -4383.545
$ws.Range("N89").Value = -24647.8335

$ws.Range("H105").Value = 3430.6924
$ws.Range("I105").Value = 2832.889
$ws.Range("J105").Value = 4775.75
$ws.Range("K105").Value = 2832.889
$ws.Range("L105").Value = 4775.75
$ws.Range("M105").Value = -1085.889
$ws.Range("N105").Value = -8269.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 0
$ws.Range("J13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("N13").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5887.7144
$ws.Range("I7").Value = 5939.9165
$ws.Range("J7").Value = 5574.5
$ws.Range("K7").Value = 17819.7495
$ws.Range("L7").Value = 16723.5
$ws.Range("M7").Value = -17707.7495
$ws.Range("N7").Value = -16947.5

$ws.Range("H8").Value = 166666940
$ws.Range("I8").Value = 166666940
$ws.Range("K8").Value = 500000820
$ws.Range("M8").Value = -500000681

$ws.Range("H22").Value = 666.3333
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 999
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 2997
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -3335

$ws.Range("H27").Value = 666.3333
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 999
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 2997
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -3201

$ws.Range("H32").Value = 9413.5
$ws.Range("I32").Value = 8883
$ws.Range("K32").Value = 26649
$ws.Range("M32").Value = -26366

$ws.Range("H107").Value = 359
$ws.Range("J107").Value = 337.06668
$ws.Range("L107").Value = 1011.20004
$ws.Range("N107").Value = -4851.20004

$ws.Range("H131").Value = 2394.862
$ws.Range("I131").Value = 1399.375
$ws.Range("J131").Value = 2774.0952
$ws.Range("K131").Value = 4198.125
$ws.Range("L131").Value = 8322.285600000001
$ws.Range("M131").Value = 841.875
$ws.Range("N131").Value = -18402.2856

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 872.5
$ws.Range("I107").Value = 691.5
$ws.Range("K107").Value = 691.5
$ws.Range("M107").Value = 1228.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2839.4285
$ws.Range("I22").Value = 1175.4
$ws.Range("J22").Value = 6999.5
$ws.Range("K22").Value = 1175.4
$ws.Range("L22").Value = 6999.5
$ws.Range("M22").Value = -880.4000000000001
$ws.Range("N22").Value = -7589.5

$ws.Range("H27").Value = 2839.4285
$ws.Range("I27").Value = 1175.4
$ws.Range("J27").Value = 6999.5
$ws.Range("K27").Value = 1175.4
$ws.Range("L27").Value = 6999.5
$ws.Range("M27").Value = -1068.4
$ws.Range("N27").Value = -7213.5

$ws.Range("H61").Value = 1536.6666
$ws.Range("I61").Value = 1536.6666
$ws.Range("K61").Value = 1536.6666
$ws.Range("M61").Value = -1334.6666

$ws.Range("H113").Value = 1536.6666
$ws.Range("I113").Value = 1536.6666
$ws.Range("K113").Value = 1536.6666
$ws.Range("M113").Value = 633.3334

$ws.Range("H132").Value = 10596.667
$ws.Range("I132").Value = 11305.454
$ws.Range("K132").Value = 33916.362
$ws.Range("M132").Value = -31386.362

$ws.Range("H136").Value = 3170.6191
$ws.Range("I136").Value = 3965.6667
$ws.Range("J136").Value = 2574.3333
$ws.Range("K136").Value = 11897.0001
$ws.Range("L136").Value = 7722.999899999999
$ws.Range("M136").Value = -9347.000100000001
$ws.Range("N136").Value = -12822.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 221.7
$ws.Range("I107").Value = 221.7
$ws.Range("K107").Value = 665.0999999999999
$ws.Range("M107").Value = 1254.9

$ws.Range("H122").Value = 4108.7144
$ws.Range("I122").Value = 3526.111
$ws.Range("J122").Value = 5157.4
$ws.Range("K122").Value = 10578.333
$ws.Range("L122").Value = 15472.2
$ws.Range("M122").Value = -8128.332999999999
$ws.Range("N122").Value = -20372.2

$ws.Range("H132").Value = 9097.895
$ws.Range("I132").Value = 5738.75
$ws.Range("J132").Value = 14856.429
$ws.Range("K132").Value = 17216.25
$ws.Range("L132").Value = 44569.287
$ws.Range("M132").Value = -14686.25
$ws.Range("N132").Value = -49629.287

$ws.Range("H136").Value = 3670.2058
$ws.Range("I136").Value = 3489.724
$ws.Range("J136").Value = 4717
$ws.Range("K136").Value = 10469.172
$ws.Range("L136").Value = 14151
$ws.Range("M136").Value = -7919.172
$ws.Range("N136").Value = -19251
